$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New "Dátum" column header (column D)
$ws.Range("D1").Value = "Dátum"
$ws.Range("D1").HorizontalAlignment = -4108   # xlCenter

# Date values for column D (Excel 1900-epoch serials, 2005-10-14 .. 2005-10-23)
$dates = 38639,38640,38641,38642,38643,38644,38645,38646,38647,38648
for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $dates[$i]
}
$ws.Range("D2:D11").NumberFormat = "mm-dd-yy"

# Explicit number formats for the existing columns (picked up by the
# original "Format Cells" pass that also introduced the date column)
$ws.Range("C2:C11").NumberFormat = "0.00"
$ws.Range("A2:A11").NumberFormat = "0"
$ws.Range("B2:B11").NumberFormat = "@"

# Column D width (auto-fit to content, like the other columns)
$ws.Columns.Item(4).AutoFit()

$ws.Range("F3").Select()
